# Fixed a bug in calcWaysRTP2
# The data rows (2-21) on Sheet1 need to be reordered to reflect the
# corrected calculation order. Row 22 (and the totals row) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(801, 3, 67, 65, 52, 45),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(301, 6, 45, 30, 60, 45),
    @(401, 9, 48, 67, 75, 45),
    @(601, 9, 60, 67, 60, 42),
    @(101, 9, 30, 15, 60, 15),
    @(501, 9, 52, 30, 75, 45),
    @(201, 9, 30, 15, 45, 30),
    @(1001, 18, 30, 75, 60, 72),
    @(701, 3, 90, 45, 97, 15),
    @(1203, 3, 15, 15, 15, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(802, 0, 4, 5, 4, 0),
    @(1, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
